$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell $ws "D2" "81.739.44"
Set-TextCell $ws "E2" "  +2.47%  "
Set-TextCell $ws "D3" "3.155.41"
Set-TextCell $ws "E3" "  -1.47%  "
Set-TextCell $ws "E4" "  +0.06%  "
Set-TextCell $ws "D5" "215.75"
Set-TextCell $ws "E5" "  +4.82%  "
Set-TextCell $ws "D6" "616.61"
Set-TextCell $ws "E6" "  -3.19%  "
Set-TextCell $ws "D7" "0.286"
Set-TextCell $ws "E7" "  +16.05%  "
Set-TextCell $ws "D8" "1.00"
Set-TextCell $ws "E8" "  -0.01%  "
Set-TextCell $ws "D9" "0.580"
Set-TextCell $ws "E9" "  -4.91%  "
Set-TextCell $ws "D10" "3.155.62"
Set-TextCell $ws "E10" "  -1.44%  "
Set-TextCell $ws "D11" "0.597"
Set-TextCell $ws "E11" "  -4.02%  "
Set-TextCell $ws "D12" "0.0000255"
Set-TextCell $ws "E12" "  +0.27%  "
Set-TextCell $ws "E13" "  -0.58%  "
Set-TextCell $ws "E14" "  -2.60%  "
Set-TextCell $ws "D15" "3.744.98"
Set-TextCell $ws "E15" "  -1.32%  "
Set-TextCell $ws "D16" "31.90"
Set-TextCell $ws "E16" "  -1.78%  "
Set-TextCell $ws "D17" "81.637.90"
Set-TextCell $ws "E17" "  +2.57%  "
Set-TextCell $ws "D18" "3.161.04"
Set-TextCell $ws "E18" "  -1.15%  "
Set-TextCell $ws "D19" "3.18"
Set-TextCell $ws "E19" "  +7.91%  "
Set-TextCell $ws "D20" "13.95"
Set-TextCell $ws "E20" "  -4.83%  "
Set-TextCell $ws "D21" "432.93"
Set-TextCell $ws "E21" "  -2.27%  "
Set-TextCell $ws "D22" "8.88"
Set-TextCell $ws "E22" "  -6.14%  "
Set-TextCell $ws "E23" "  -2.99%  "
Set-TextCell $ws "D24" "7.25"
Set-TextCell $ws "E24" "  +5.42%  "
Set-TextCell $ws "D25" "5.22"
Set-TextCell $ws "E25" "  +7.31%  "
Set-TextCell $ws "D26" "11.66"
Set-TextCell $ws "E26" "  +6.90%  "
Set-TextCell $ws "D27" "3.335.23"
Set-TextCell $ws "E27" "  -1.01%  "
Set-TextCell $ws "D28" "76.40"
Set-TextCell $ws "E28" "  -1.83%  "
Set-TextCell $ws "E29" "  -0.42%  "
Set-TextCell $ws "D30" "0.0000120"
Set-TextCell $ws "E30" "  +0.79%  "
Set-TextCell $ws "D31" "1.00"
Set-TextCell $ws "E31" "  +0.01%  "
Set-TextCell $ws "D32" "8.97"
Set-TextCell $ws "E32" "  -3.26%  "
Set-TextCell $ws "D33" "564.78"
Set-TextCell $ws "E33" "  +5.73%  "
Set-TextCell $ws "D34" "1.48"
Set-TextCell $ws "E34" "  -1.06%  "
Set-TextCell $ws "E35" "  +18.34%  "
Set-TextCell $ws "D36" "0.151"
Set-TextCell $ws "E36" "  +2.74%  "
Set-TextCell $ws "E37" "  -2.01%  "
Set-TextCell $ws "D38" "22.53"
Set-TextCell $ws "E38" "  -3.50%  "
Set-TextCell $ws "E39" "  -0.05%  "
Set-TextCell $ws "D40" "6.06"
Set-TextCell $ws "E40" "  +8.62%  "
Set-TextCell $ws "D41" "0.403"
Set-TextCell $ws "E41" "  -1.75%  "
Set-TextCell $ws "D42" "20.82"
Set-TextCell $ws "E42" "  +3.87%  "
Set-TextCell $ws "D43" "3.00"
Set-TextCell $ws "E43" "  +15.37%  "
Set-TextCell $ws "D44" "1.99"
Set-TextCell $ws "E44" "  +8.49%  "
Set-TextCell $ws "D45" "159.80"
Set-TextCell $ws "E45" "  -3.05%  "
Set-TextCell $ws "E46" "  +0.03%  "
Set-TextCell $ws "D47" "186.07"
Set-TextCell $ws "E47" "  -3.35%  "
Set-TextCell $ws "D48" "44.44"
Set-TextCell $ws "E48" "  +1.85%  "
Set-TextCell $ws "E49" "  -0.43%  "
Set-TextCell $ws "D50" "26.31"
Set-TextCell $ws "E50" "  +2.07%  "
Set-TextCell $ws "D51" "0.760"
Set-TextCell $ws "E51" "  -6.05%  "
